$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '28.932.81'
$ws.Range('E2').Value = '  +1.72%  '
Set-TextValue $ws.Range('D3') '1.893.24'
$ws.Range('E3').Value = '  +1.57%  '
Set-TextValue $ws.Range('D4') '1.002'
$ws.Range('E4').Value = '  -0.47%  '
Set-TextValue $ws.Range('D5') '326.14'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('E7').Value = '  +0.86%  '
Set-TextValue $ws.Range('D8') '0.3906'
$ws.Range('E8').Value = '  +2.06%  '
Set-TextValue $ws.Range('D9') '0.07847'
$ws.Range('E9').Value = '  +0.46%  '
Set-TextValue $ws.Range('D10') '0.9895'
$ws.Range('E10').Value = '  +0.19%  '
Set-TextValue $ws.Range('D11') '21.90'
$ws.Range('E11').Value = '  +1.66%  '
Set-TextValue $ws.Range('D12') '1.885.41'
$ws.Range('E12').Value = '  +0.34%  '
Set-TextValue $ws.Range('D13') '7.041'
$ws.Range('E13').Value = '  +2.15%  '
Set-TextValue $ws.Range('D14') '5.706'
$ws.Range('E14').Value = '  +1.53%  '
Set-TextValue $ws.Range('D15') '0.06959'
$ws.Range('E15').Value = '  +1.01%  '
Set-TextValue $ws.Range('D16') '88.08'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('E17').Value = '  -0.27%  '
Set-TextValue $ws.Range('D18') '0.000009959'
$ws.Range('E18').Value = '  +0.28%  '
Set-TextValue $ws.Range('D19') '16.93'
$ws.Range('E19').Value = '  +1.78%  '
Set-TextValue $ws.Range('D20') '1.002'
$ws.Range('E20').Value = '  -0.26%  '
Set-TextValue $ws.Range('D21') '28.925.04'
$ws.Range('E21').Value = '  +1.64%  '
Set-TextValue $ws.Range('D22') '5.292'
$ws.Range('E22').Value = '  +0.95%  '
Set-TextValue $ws.Range('D23') '11.00'
$ws.Range('E23').Value = '  +1.15%  '
Set-TextValue $ws.Range('D24') '2.132.17'
$ws.Range('E24').Value = '  +1.62%  '
Set-TextValue $ws.Range('D25') '2.061'
$ws.Range('E25').Value = '  -1.86%  '
Set-TextValue $ws.Range('D26') '156.06'
Set-TextValue $ws.Range('D27') '19.31'
$ws.Range('E27').Value = '  +1.16%  '
Set-TextValue $ws.Range('D28') '5.928'
$ws.Range('E28').Value = '  +4.65%  '
Set-TextValue $ws.Range('D29') '1.932'
$ws.Range('E29').Value = '  +0.92%  '
Set-TextValue $ws.Range('D30') '117.68'
$ws.Range('E30').Value = '  +0.27%  '
Set-TextValue $ws.Range('D31') '0.09346'
$ws.Range('E31').Value = '  +0.84%  '
Set-TextValue $ws.Range('D32') '0.9112'
$ws.Range('E32').Value = '  +0.71%  '
Set-TextValue $ws.Range('D33') '5.296'
$ws.Range('E33').Value = '  +0.74%  '
Set-TextValue $ws.Range('D34') '1.336'
$ws.Range('E34').Value = '  +1.38%  '
Set-TextValue $ws.Range('D35') '3.278'
$ws.Range('E35').Value = '  -0.51%  '
Set-TextValue $ws.Range('D36') '1.191'
$ws.Range('E36').Value = '  +4.05%  '
$ws.Range('E37').Value = '  +1.35%  '
Set-TextValue $ws.Range('D38') '0.02074'
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('E39').Value = '  -0.32%  '
Set-TextValue $ws.Range('D40') '7.766'
$ws.Range('E40').Value = '  +1.44%  '
Set-TextValue $ws.Range('D41') '0.5682'
$ws.Range('E41').Value = '  +2.40%  '
$ws.Range('E42').Value = '  +0.28%  '
Set-TextValue $ws.Range('D43') '9.780'
$ws.Range('E43').Value = '  +1.60%  '
Set-TextValue $ws.Range('D44') '2.251'
$ws.Range('E44').Value = '  +6.21%  '
Set-TextValue $ws.Range('D45') '11.90'
$ws.Range('E45').Value = '  +3.34%  '
Set-TextValue $ws.Range('D46') '0.5355'
$ws.Range('E46').Value = '  +2.42%  '
$ws.Range('E47').Value = '  -0.72%  '
Set-TextValue $ws.Range('D48') '1.842'
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('E49').Value = '  +4.36%  '
Set-TextValue $ws.Range('D50') '112.42'
$ws.Range('E50').Value = '  +0.47%  '
Set-TextValue $ws.Range('D51') '1.069'
$ws.Range('E51').Value = '  -5.29%  '
